$d = $word.ActiveDocument

# Remove the literal "<lb/>" marker run that immediately follows
# "... enfin empescheroit quelle ne" (the text is unique in the document,
# so this precisely targets the run to delete without touching any of
# the many other "<lb/>" markers elsewhere in the transcription).
$found = $d.Content.Find.Execute(
    "quelle ne<lb/>", $false, $false, $false, $false, $false,
    $true, 1, $false, "quelle ne", 2)

if (-not $found) {
    throw "Could not find the target '<lb/>' marker to remove"
}

Write-Output "Replaced: $found"
